$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A182").Value = "2023-12-11 12:01:35"
$ws.Range("B182").Value = 0.0004

$ws.Range("A183").Value = "2023-12-11 12:01:47"
$ws.Range("B183").Value = 0.0006000000000000001

$ws.Range("A184").Value = "2023-12-11 12:02:04"
$ws.Range("B184").Value = 0.0008

$ws.Range("A185").Value = "2023-12-11 12:02:13"
$ws.Range("B185").Value = 0.0002
